$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("J3").Value = 6657
$ws.Range("N3").Value = -6885
$ws.Range("H3").Value = 6657
$ws.Range("L3").Value = 6657
# Row 19
$ws.Range("L19").Value = 3080.5
$ws.Range("J19").Value = 3080.5
$ws.Range("I19").Value = 512
$ws.Range("N19").Value = -3430.5
$ws.Range("H19").Value = 1653.5555
$ws.Range("M19").Value = -337
$ws.Range("K19").Value = 512
# Row 76
$ws.Range("L76").Value = 142860400
$ws.Range("H76").Value = 66670110
$ws.Range("J76").Value = 142860400
$ws.Range("N76").Value = -142861030
# Row 79
$ws.Range("H79").Value = 66670110
$ws.Range("J79").Value = 142860400
$ws.Range("N79").Value = -142862584
$ws.Range("L79").Value = 142860400
# Row 98
$ws.Range("J98").Value = 3849
$ws.Range("I98").Value = 2607.7778
$ws.Range("L98").Value = 3849
$ws.Range("N98").Value = -6845
$ws.Range("K98").Value = 2607.7778
$ws.Range("M98").Value = -1109.7778
$ws.Range("H98").Value = 2989.6924
# Row 102
$ws.Range("J102").Value = 6657
$ws.Range("H102").Value = 6657
$ws.Range("N102").Value = -13147
$ws.Range("L102").Value = 6657
# Row 112
$ws.Range("H112").Value = 2099.4482
$ws.Range("L112").Value = 6298.344599999999
$ws.Range("J112").Value = 2099.4482
$ws.Range("N112").Value = -8514.3446
# Row 113
$ws.Range("H113").Value = 3422.2222
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("J113").Value = 4000
$ws.Range("M113").Value = 554
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -10508
# Row 122
$ws.Range("N122").Value = -16447
$ws.Range("J122").Value = 3849
$ws.Range("L122").Value = 11547
$ws.Range("I122").Value = 2607.7778
$ws.Range("H122").Value = 2989.6924
$ws.Range("M122").Value = -5373.3334
$ws.Range("K122").Value = 7823.3334
# Row 125
$ws.Range("M125").ClearContents() | Out-Null
$ws.Range("K125").Value = 10800
$ws.Range("N125").Value = -8340
$ws.Range("H125").Value = 1200
$ws.Range("I125").Value = 1200
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
# Row 127
$ws.Range("J127").Value = 1280.75
$ws.Range("K127").Value = 963.4999799999999
$ws.Range("I127").Value = 321.16666
$ws.Range("N127").Value = -13762.25
$ws.Range("M127").Value = 3996.50002
$ws.Range("H127").Value = 561.0625
$ws.Range("L127").Value = 3842.25
# Row 137
$ws.Range("N137").Value = -9720.375
$ws.Range("J137").Value = 1540.125
$ws.Range("M137").Value = -1623.4998
$ws.Range("I137").Value = 1391.1666
$ws.Range("K137").Value = 4173.4998
$ws.Range("H137").Value = 1450.75
$ws.Range("L137").Value = 4620.375
# Row 141
$ws.Range("M141").Value = 2504.6
$ws.Range("H141").Value = 891.8
$ws.Range("K141").Value = 2675.4
$ws.Range("I141").Value = 891.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 9192.75
$ws.Range("I2").Value = 811.1111
$ws.Range("K2").Value = 811.1111
$ws.Range("M2").Value = -698.1111
# Row 32
$ws.Range("H32").Value = 3071.5293
$ws.Range("K32").Value = 2679.8958
$ws.Range("I32").Value = 2679.8958
$ws.Range("M32").Value = -2392.8958
# Row 61
$ws.Range("H61").Value = 2066.6667
$ws.Range("I61").Value = 0
$ws.Range("M61").ClearContents() | Out-Null
$ws.Range("K61").Value = 0
# Row 101
$ws.Range("L101").ClearContents() | Out-Null
$ws.Range("H101").Value = 0
$ws.Range("N101").Value = 0
$ws.Range("J101").Value = 0
# Row 116
$ws.Range("K116").Value = 811.1111
$ws.Range("H116").Value = 9192.75
$ws.Range("M116").Value = 1482.8889
$ws.Range("I116").Value = 811.1111
# Row 122
$ws.Range("I122").Value = 4000
$ws.Range("H122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550
# Row 132
$ws.Range("J132").Value = 2790.75
$ws.Range("L132").Value = 8372.25
$ws.Range("H132").Value = 2244
$ws.Range("N132").Value = -13432.25
# Row 136
$ws.Range("K136").Value = 0
$ws.Range("H136").Value = 2066.6667
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("I136").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("K3").Value = 811.1111
$ws.Range("M3").Value = -697.1111
$ws.Range("I3").Value = 811.1111
$ws.Range("H3").Value = 9192.75
# Row 128
$ws.Range("M128").ClearContents() | Out-Null
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("H128").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("M99").Value = -357.7
$ws.Range("K99").Value = 1855.7
$ws.Range("H99").Value = 1868.8182
$ws.Range("I99").Value = 1855.7
# Row 126
$ws.Range("K126").Value = 5567.1
$ws.Range("M126").Value = -3097.1
$ws.Range("H126").Value = 1868.8182
$ws.Range("I126").Value = 1855.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("K5").Value = 1675.5
$ws.Range("H5").Value = 558.5
$ws.Range("I5").Value = 558.5
$ws.Range("M5").Value = -1563.5
# Row 113
$ws.Range("H113").Value = 610.9394
$ws.Range("L113").Value = 1973.4444
$ws.Range("J113").Value = 657.8148
$ws.Range("N113").Value = -6313.4444
# Row 116
$ws.Range("H116").Value = 2869.8
$ws.Range("J116").Value = 2712.5
$ws.Range("L116").Value = 8137.5
$ws.Range("N116").Value = -15021.5
# Row 122
$ws.Range("I122").Value = 499.2
$ws.Range("H122").Value = 463.83334
$ws.Range("K122").Value = 4492.8
$ws.Range("M122").Value = -2042.8
# Row 135
$ws.Range("H135").Value = 558.5
$ws.Range("I135").Value = 558.5
$ws.Range("M135").Value = -2491.5
$ws.Range("K135").Value = 5026.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1584.8667
$ws.Range("I113").Value = 1035.7778
$ws.Range("K113").Value = 1035.7778
$ws.Range("J113").Value = 2408.5
$ws.Range("M113").Value = 1134.2222
$ws.Range("L113").Value = 2408.5
$ws.Range("N113").Value = -6748.5
# Row 122
$ws.Range("N122").Value = -8500
$ws.Range("J122").Value = 1200
$ws.Range("L122").Value = 3600
$ws.Range("I122").Value = 1007
$ws.Range("H122").Value = 1103.5
$ws.Range("M122").Value = -571
$ws.Range("K122").Value = 3021
# Row 126
$ws.Range("K126").Value = 5377.2498
$ws.Range("H126").Value = 1957.625
$ws.Range("M126").Value = -2907.2498
$ws.Range("N126").Value = -12299.75
$ws.Range("L126").Value = 7359.75
$ws.Range("J126").Value = 2453.25
$ws.Range("I126").Value = 1792.4166

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("L7").Value = 3301.25
$ws.Range("M7").Value = -2688
$ws.Range("I7").Value = 2800
$ws.Range("J7").Value = 3301.25
$ws.Range("N7").Value = -3525.25
$ws.Range("K7").Value = 2800
$ws.Range("H7").Value = 3000.5
# Row 40
$ws.Range("K40").Value = 2896.375
$ws.Range("I40").Value = 2896.375
$ws.Range("M40").Value = -2760.375
$ws.Range("H40").Value = 2896.375
# Row 126
$ws.Range("K126").Value = 8400
$ws.Range("H126").Value = 3000.5
$ws.Range("M126").Value = -5930
$ws.Range("N126").Value = -14843.75
$ws.Range("L126").Value = 9903.75
$ws.Range("J126").Value = 3301.25
$ws.Range("I126").Value = 2800
# Row 132
$ws.Range("L132").Value = 136241.34
$ws.Range("K132").Value = 3987.6426
$ws.Range("J132").Value = 45413.78
$ws.Range("I132").Value = 1329.2142
$ws.Range("H132").Value = 21210.49
$ws.Range("N132").Value = -141301.34
$ws.Range("M132").Value = -1457.6426

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("L101").Value = 14840.4
$ws.Range("H101").Value = 14840.4
$ws.Range("N101").Value = -21330.4
$ws.Range("J101").Value = 14840.4
# Row 126
$ws.Range("K126").Value = 187502364
$ws.Range("H126").Value = 62500788
$ws.Range("M126").ClearContents() | Out-Null
$ws.Range("N126").Value = -187499894
$ws.Range("L126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("I126").Value = 62500788
# Row 132
$ws.Range("L132").Value = 12696.4995
$ws.Range("K132").Value = 9105.8181
$ws.Range("J132").Value = 4232.1665
$ws.Range("I132").Value = 3035.2727
$ws.Range("H132").Value = 3457.7058
$ws.Range("N132").Value = -17756.4995
$ws.Range("M132").Value = -6575.8181
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("N136").Value = 0
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("I136").Value = 0
$ws.Range("L136").ClearContents() | Out-Null
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
